$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Price (D) and Volume/Change-1h (E) columns row by row with the
# latest scraped cryptocurrency figures.
#
# The Price column stores numeric-looking text (e.g. "511.07") as plain text
# in the source workbook. To stop Excel auto-converting such strings into real
# numbers (which would also introduce floating point artifacts, e.g.
# 511.06999999999999) we briefly force a Text number format before writing the
# value, then restore the default "Normal" style so the cell formatting is
# left exactly as it was.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.620.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.76%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.030.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.26%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "511.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.19%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.440"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.48"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.110"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.97%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.366"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.88%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.555.53"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.40%  "

$ws.Range("E13").Value = "  +2.27%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.60%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000165"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +10.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "57.682.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.76%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +9.33%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.035.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.98"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "333.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.77"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.97%  "

$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.497"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.61%  "

$ws.Range("E26").Value = "  +3.84%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.17%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0923"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.79"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.66%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.41"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +10.50%  "

$ws.Range("E31").Value = "  +4.11%  "

$ws.Range("E32").Value = "  +2.88%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.18%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "156.20"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.70"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.27%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.12%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.77%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.63"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.96%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0684"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.068.68"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.39%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.48"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.72%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.87"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.60%  "

$ws.Range("E43").Value = "  +0.08%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.658"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.33%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.300.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.93%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.68%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.987"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.27%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.04%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0239"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.84%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.44"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.38%  "

$ws.Range("E51").Value = "  -3.71%  "
